$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> @{ D = new price string (optional); E = new volume string }
$updates = @{
    2  = @{ D = "64.108.42";  E = "  -0.28%  " }
    3  = @{ D = "3.394.00";   E = "  -1.85%  " }
    4  = @{            E = "  +0.07%  " }
    5  = @{ D = "568.58";     E = "  -1.23%  " }
    6  = @{ D = "155.93";     E = "  -2.96%  " }
    7  = @{ D = "0.612";      E = "  +4.54%  " }
    8  = @{            E = "  +0.12%  " }
    9  = @{ D = "3.396.64";   E = "  -1.71%  " }
    10 = @{ D = "7.17";       E = "  -2.06%  " }
    11 = @{            E = "  -3.93%  " }
    12 = @{ D = "0.439";      E = "  -1.51%  " }
    13 = @{ D = "3.988.61";   E = "  -1.65%  " }
    14 = @{            E = "  -0.23%  " }
    15 = @{ D = "0.0000186";  E = "  -4.20%  " }
    16 = @{ D = "27.71";      E = "  -4.91%  " }
    17 = @{ D = "64.244.18";  E = "  -0.09%  " }
    18 = @{ D = "3.412.45";   E = "  -1.19%  " }
    19 = @{            E = "  -1.70%  " }
    20 = @{ D = "13.92";      E = "  -4.01%  " }
    21 = @{ D = "373.72";     E = "  -3.68%  " }
    22 = @{ D = "7.92";       E = "  -3.98%  " }
    23 = @{ D = "0.546";      E = "  -0.58%  " }
    24 = @{ D = "0.997";      E = "  -0.42%  " }
    25 = @{ D = "71.65";      E = "  -2.40%  " }
    26 = @{ D = "0.0000118";  E = "  -5.68%  " }
    27 = @{ D = "9.93";       E = "  +4.09%  " }
    28 = @{ D = "0.176";      E = "  -2.05%  " }
    29 = @{ D = "0.999";      E = "  -0.38%  " }
    30 = @{ D = "1.47";       E = "  +1.42%  " }
    31 = @{ D = "6.06";       E = "  -1.69%  " }
    32 = @{            E = "  -1.28%  " }
    33 = @{ D = "23.07";      E = "  -2.89%  " }
    34 = @{ D = "7.19";       E = "  +0.69%  " }
    35 = @{            E = "  +4.78%  " }
    36 = @{ D = "159.62";     E = "  -0.74%  " }
    37 = @{            E = "  -0.44%  " }
    38 = @{ D = "0.0758";     E = "  -2.68%  " }
    39 = @{ D = "26.71";      E = "  -2.91%  " }
    40 = @{ D = "6.71";       E = "  +1.71%  " }
    41 = @{ D = "2.825.36";   E = "  -3.32%  " }
    42 = @{ D = "4.59";       E = "  +1.54%  " }
    44 = @{ D = "0.0311";     E = "  -3.46%  " }
    45 = @{ D = "0.762";      E = "  -0.99%  " }
    46 = @{            E = "  +6.52%  " }
    47 = @{ D = "1.07";       E = "  -1.76%  " }
    48 = @{ D = "312.19";     E = "  +5.38%  " }
    49 = @{ D = "0.109";      E = "  +0.54%  " }
    50 = @{ D = "6.54";       E = "  -1.13%  " }
    51 = @{ D = "0.854";      E = "  -1.63%  " }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($vals.ContainsKey("D")) {
        $cell = $ws.Range("D$row")
        $cell.NumberFormat = "@"
        $cell.Value = $vals["D"]
    }
    if ($vals.ContainsKey("E")) {
        $cell = $ws.Range("E$row")
        $cell.NumberFormat = "@"
        $cell.Value = $vals["E"]
    }
}
